# [Feat] Setup DB from scraping
# Adds a "naverPlaceUrl" column (AF) with per-restaurant Naver Place
# short-links (as real hyperlinks on most rows), and removes the stray
# 26th row that held a leftover/duplicate kakao review URL fragment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the orphan row 26 (AB26/AC26/AD26/AE26 leftover scrap data).
#    This also lets the engine garbage-collect the now-unused shared
#    string "https://place.map.kakao.com/25462155" on save.
# ---------------------------------------------------------------------
$ws.Rows.Item(26).Delete()

# ---------------------------------------------------------------------
# 2) New column AF: "naverPlaceUrl"
# ---------------------------------------------------------------------
$ws.Range("AF1").Value = "naverPlaceUrl"

# row -> Naver Place short URL (rows 5 & 23 intentionally have no URL)
$naverUrls = [ordered]@{
    2  = "https://naver.me/5GywN2l6"
    3  = "https://naver.me/5On93MHs"
    4  = "https://naver.me/x9LgiJ2R"
    6  = "https://naver.me/GNy6CHfW"
    7  = "https://naver.me/xmrHAKBa"
    8  = "https://naver.me/5L3dgXyA"
    9  = "https://naver.me/GQ4OFWHF"
    10 = "https://naver.me/xaPuSkys"
    11 = "https://naver.me/5X9p62AZ"
    12 = "https://naver.me/GkKhCpEv"
    13 = "https://naver.me/5zJkKdG4"
    14 = "https://naver.me/5DHgpYCH"
    15 = "https://naver.me/FCbsFDHb"
    16 = "https://naver.me/GEAnmV7k"
    17 = "https://naver.me/FlJwqsOD"
    18 = "https://naver.me/FKKG0F00"
    19 = "https://naver.me/GcjKdEIz"
    20 = "https://naver.me/xrPzfxWi"
    21 = "https://naver.me/5AmSDHn6"
    22 = "https://naver.me/59jXODvJ"
    24 = "https://naver.me/FPsngYjo"
    25 = "https://naver.me/I5FvnEVf"
}

# Row 6 only gets plain text (no hyperlink, no hyperlink style) -- matches
# source data where that particular link wasn't wired up as a clickable
# hyperlink like its neighbours.
$ws.Range("AF6").Value = $naverUrls[6]

# Row 5 has no URL at all yet, but already carries the hyperlink style
# (formatting got pre-applied down the column before the data landed).
$ws.Range("AC2").Copy()
$ws.Range("AF5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows that get real, clickable hyperlinks -- added in this exact order
# so the relationship ids (rId3, rId4, ...) line up with the cell order.
$hyperlinkRows = @(2, 3, 4, 7, 8, 9, 10, 11, 12, 13, 15, 14, 16, 17, 18, 19, 20, 21, 22, 24, 25)
foreach ($r in $hyperlinkRows) {
    $addr = "AF" + $r
    $ws.Hyperlinks.Add($ws.Range($addr), $naverUrls[$r])
}

# Re-apply the workbook's existing hyperlink style (xfId 10 / "Hyperlink")
# to each newly-linked cell so it matches the look of AC2/AC11.
$ws.Range("AC2").Copy()
foreach ($r in $hyperlinkRows) {
    $addr = "AF" + $r
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Final selection matches where the editor ended up after the edit.
# ---------------------------------------------------------------------
$ws.Range("AF25").Select()
